$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 7 more "Named" method sample rows (rows 6-12) below the existing
# data (rows 1-5), same layout: Date, Method, ElapsedMs, wordCount,
# sentenceCount, posWordCount, negWordCount, posWordPercentage,
# negWordPercentage, positivePhraseCount, negativePhraseCount,
# posPhrasePercentage, negPhrasePercentage.
$rows = @(
    @(6,  42601.767372685186, "Named", 3252, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(7,  42601.769293981481, "Named", 3159, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(8,  42601.770046296297, "Named", 3107, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(9,  42601.771041666667, "Named", 3051, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(10, 42601.771458333336, "Named", 3199, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(11, 42601.772812499999, "Named", 3110, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(12, 42601.773298611108, "Named", 3306, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
}

# Column A widened slightly to keep fitting the date/time values.
$ws.Columns.Item(1).ColumnWidth = 14
